# Appends the 8 simulation result rows (295-302) recorded for the new
# injection-rate / tx-per-block sweep onto the existing results table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{
        Row = 295
        Cells = @(
            @{ Col = 1; Type = 'n'; Value = 2 },
            @{ Col = 2; Type = 'n'; Value = 1 },
            @{ Col = 3; Type = 'n'; Value = 2 },
            @{ Col = 4; Type = 'n'; Value = 10 },
            @{ Col = 5; Type = 'n'; Value = 10 },
            @{ Col = 6; Type = 'n'; Value = 3 },
            @{ Col = 7; Type = 'n'; Value = 50 },
            @{ Col = 8; Type = 'n'; Value = 0 },
            @{ Col = 9; Type = 'n'; Value = 0 },
            @{ Col = 10; Type = 'n'; Value = 100 },
            @{ Col = 11; Type = 'n'; Value = 200 },
            @{ Col = 12; Type = 'b'; Value = $false },
            @{ Col = 13; Type = 'n'; Value = -2 },
            @{ Col = 14; Type = 'n'; Value = 100 },
            @{ Col = 15; Type = 'n'; Value = 256 },
            @{ Col = 16; Type = 'n'; Value = 10 },
            @{ Col = 17; Type = 's'; Value = "<-parameter / result->" },
            @{ Col = 18; Type = 'n'; Value = 8 },
            @{ Col = 19; Type = 'n'; Value = -1; Fmt = '0.000000' },
            @{ Col = 20; Type = 'n'; Value = 40.37446797688803; Fmt = '0.000000' },
            @{ Col = 21; Type = 'n'; Value = 13.95394086837769; Fmt = '0.000000' },
            @{ Col = 22; Type = 'n'; Value = 0.1116315269470215; Fmt = '0.000000' },
            @{ Col = 23; Type = 'n'; Value = [double]"8.034706115722656e-05"; Fmt = '0.000000' },
            @{ Col = 24; Type = 'n'; Value = 0.0725; Fmt = '0.000000' },
            @{ Col = 25; Type = 'n'; Value = 0.03625; Fmt = '0.000000' },
            @{ Col = 26; Type = 'n'; Value = 0; Fmt = '0.000000' },
            @{ Col = 27; Type = 'n'; Value = 1.221457004547119; Fmt = '0.000' },
            @{ Col = 28; Type = 'n'; Value = 742.4; Fmt = '0.00' },
            @{ Col = 29; Type = 'n'; Value = 0; Fmt = '0.00' }
        )
    },
    @{
        Row = 296
        Cells = @(
            @{ Col = 1; Type = 'n'; Value = 2 },
            @{ Col = 2; Type = 'n'; Value = 1 },
            @{ Col = 3; Type = 'n'; Value = 2 },
            @{ Col = 4; Type = 'n'; Value = 10 },
            @{ Col = 5; Type = 'n'; Value = 10 },
            @{ Col = 6; Type = 'n'; Value = 3 },
            @{ Col = 7; Type = 'n'; Value = 50 },
            @{ Col = 8; Type = 'n'; Value = 0 },
            @{ Col = 9; Type = 'n'; Value = 0 },
            @{ Col = 10; Type = 'n'; Value = 100 },
            @{ Col = 11; Type = 'n'; Value = 200 },
            @{ Col = 12; Type = 'b'; Value = $false },
            @{ Col = 13; Type = 'n'; Value = -2 },
            @{ Col = 14; Type = 'n'; Value = 100 },
            @{ Col = 15; Type = 'n'; Value = 600 },
            @{ Col = 16; Type = 'n'; Value = 20 },
            @{ Col = 17; Type = 's'; Value = "<-parameter / result->" },
            @{ Col = 18; Type = 'n'; Value = 5 },
            @{ Col = 19; Type = 'n'; Value = -1; Fmt = '0.000000' },
            @{ Col = 20; Type = 'n'; Value = 29.42369191542915; Fmt = '0.000000' },
            @{ Col = 21; Type = 'n'; Value = 17.55800247192383; Fmt = '0.000000' },
            @{ Col = 22; Type = 'n'; Value = 0.08779001235961914; Fmt = '0.000000' },
            @{ Col = 23; Type = 'n'; Value = [double]"7.033348083496094e-05"; Fmt = '0.000000' },
            @{ Col = 24; Type = 'n'; Value = 0.056640625; Fmt = '0.000000' },
            @{ Col = 25; Type = 'n'; Value = 0.0283203125; Fmt = '0.000000' },
            @{ Col = 26; Type = 'n'; Value = 0; Fmt = '0.000000' },
            @{ Col = 27; Type = 'n'; Value = 1.223086357116699; Fmt = '0.000' },
            @{ Col = 28; Type = 'n'; Value = 580; Fmt = '0.00' },
            @{ Col = 29; Type = 'n'; Value = 0; Fmt = '0.00' }
        )
    },
    @{
        Row = 297
        Cells = @(
            @{ Col = 1; Type = 'n'; Value = 2 },
            @{ Col = 2; Type = 'n'; Value = 1 },
            @{ Col = 3; Type = 'n'; Value = 2 },
            @{ Col = 4; Type = 'n'; Value = 10 },
            @{ Col = 5; Type = 'n'; Value = 10 },
            @{ Col = 6; Type = 'n'; Value = 3 },
            @{ Col = 7; Type = 'n'; Value = 50 },
            @{ Col = 8; Type = 'n'; Value = 0 },
            @{ Col = 9; Type = 'n'; Value = 0 },
            @{ Col = 10; Type = 'n'; Value = 100 },
            @{ Col = 11; Type = 'n'; Value = 200 },
            @{ Col = 12; Type = 'b'; Value = $false },
            @{ Col = 13; Type = 'n'; Value = -2 },
            @{ Col = 14; Type = 'n'; Value = 100 },
            @{ Col = 15; Type = 'n'; Value = 256 },
            @{ Col = 16; Type = 'n'; Value = 10 },
            @{ Col = 17; Type = 's'; Value = "<-parameter / result->" },
            @{ Col = 18; Type = 'n'; Value = 9 },
            @{ Col = 19; Type = 'n'; Value = -1; Fmt = '0.000000' },
            @{ Col = 20; Type = 'n'; Value = 51.1702262654024; Fmt = '0.000000' },
            @{ Col = 21; Type = 'n'; Value = 17.50065485636394; Fmt = '0.000000' },
            @{ Col = 22; Type = 'n'; Value = 0.1575058937072754; Fmt = '0.000000' },
            @{ Col = 23; Type = 'n'; Value = 0.0001037120819091797; Fmt = '0.000000' },
            @{ Col = 24; Type = 'n'; Value = 0.101953125; Fmt = '0.000000' },
            @{ Col = 25; Type = 'n'; Value = 0.05097656249999999; Fmt = '0.000000' },
            @{ Col = 26; Type = 'n'; Value = 0; Fmt = '0.000000' },
            @{ Col = 27; Type = 'n'; Value = 1.225745677947998; Fmt = '0.000' },
            @{ Col = 28; Type = 'n'; Value = 1044; Fmt = '0.00' },
            @{ Col = 29; Type = 'n'; Value = 0; Fmt = '0.00' }
        )
    },
    @{
        Row = 298
        Cells = @(
            @{ Col = 1; Type = 'n'; Value = 2 },
            @{ Col = 2; Type = 'n'; Value = 1 },
            @{ Col = 3; Type = 'n'; Value = 2 },
            @{ Col = 4; Type = 'n'; Value = 10 },
            @{ Col = 5; Type = 'n'; Value = 10 },
            @{ Col = 6; Type = 'n'; Value = 3 },
            @{ Col = 7; Type = 'n'; Value = 50 },
            @{ Col = 8; Type = 'n'; Value = 0 },
            @{ Col = 9; Type = 'n'; Value = 0 },
            @{ Col = 10; Type = 'n'; Value = 100 },
            @{ Col = 11; Type = 'n'; Value = 200 },
            @{ Col = 12; Type = 'b'; Value = $false },
            @{ Col = 13; Type = 'n'; Value = -2 },
            @{ Col = 14; Type = 'n'; Value = 100 },
            @{ Col = 15; Type = 'n'; Value = 256 },
            @{ Col = 16; Type = 'n'; Value = 10 },
            @{ Col = 17; Type = 's'; Value = "<-parameter / result->" },
            @{ Col = 18; Type = 'n'; Value = 9 },
            @{ Col = 19; Type = 'n'; Value = -1; Fmt = '0.000000' },
            @{ Col = 20; Type = 'n'; Value = 50.4506430906408; Fmt = '0.000000' },
            @{ Col = 21; Type = 'n'; Value = 17.42459932963054; Fmt = '0.000000' },
            @{ Col = 22; Type = 'n'; Value = 0.1568213939666748; Fmt = '0.000000' },
            @{ Col = 23; Type = 'n'; Value = [double]"9.751319885253906e-05"; Fmt = '0.000000' },
            @{ Col = 24; Type = 'n'; Value = 0.101953125; Fmt = '0.000000' },
            @{ Col = 25; Type = 'n'; Value = 0.05097656249999999; Fmt = '0.000000' },
            @{ Col = 26; Type = 'n'; Value = 0; Fmt = '0.000000' },
            @{ Col = 27; Type = 'n'; Value = 1.224261045455933; Fmt = '0.000' },
            @{ Col = 28; Type = 'n'; Value = 1044; Fmt = '0.00' },
            @{ Col = 29; Type = 'n'; Value = 0; Fmt = '0.00' }
        )
    },
    @{
        Row = 299
        Cells = @(
            @{ Col = 1; Type = 'n'; Value = 2 },
            @{ Col = 2; Type = 'n'; Value = 1 },
            @{ Col = 3; Type = 'n'; Value = 2 },
            @{ Col = 4; Type = 'n'; Value = 10 },
            @{ Col = 5; Type = 'n'; Value = 10 },
            @{ Col = 6; Type = 'n'; Value = 3 },
            @{ Col = 7; Type = 'n'; Value = 50 },
            @{ Col = 8; Type = 'n'; Value = 0 },
            @{ Col = 9; Type = 'n'; Value = 0 },
            @{ Col = 10; Type = 'n'; Value = 100 },
            @{ Col = 11; Type = 'n'; Value = 200 },
            @{ Col = 12; Type = 'b'; Value = $false },
            @{ Col = 13; Type = 'n'; Value = -2 },
            @{ Col = 14; Type = 'n'; Value = 100 },
            @{ Col = 15; Type = 'n'; Value = 256 },
            @{ Col = 16; Type = 'n'; Value = 10 },
            @{ Col = 17; Type = 's'; Value = "<-parameter / result->" },
            @{ Col = 18; Type = 'n'; Value = 7 },
            @{ Col = 19; Type = 'n'; Value = -1; Fmt = '0.000000' },
            @{ Col = 20; Type = 'n'; Value = 32.62987295786539; Fmt = '0.000000' },
            @{ Col = 21; Type = 'n'; Value = 10.87885856628418; Fmt = '0.000000' },
            @{ Col = 22; Type = 'n'; Value = 0.07615200996398924; Fmt = '0.000000' },
            @{ Col = 23; Type = 'n'; Value = [double]"7.367134094238281e-05"; Fmt = '0.000000' },
            @{ Col = 24; Type = 'n'; Value = 0.047578125; Fmt = '0.000000' },
            @{ Col = 25; Type = 'n'; Value = 0.0237890625; Fmt = '0.000000' },
            @{ Col = 26; Type = 'n'; Value = 0; Fmt = '0.000000' },
            @{ Col = 27; Type = 'n'; Value = 1.223163843154907; Fmt = '0.000' },
            @{ Col = 28; Type = 'n'; Value = 487.2; Fmt = '0.00' },
            @{ Col = 29; Type = 'n'; Value = 0; Fmt = '0.00' }
        )
    },
    @{
        Row = 300
        Cells = @(
            @{ Col = 1; Type = 'n'; Value = 2 },
            @{ Col = 2; Type = 'n'; Value = 1 },
            @{ Col = 3; Type = 'n'; Value = 2 },
            @{ Col = 4; Type = 'n'; Value = 10 },
            @{ Col = 5; Type = 'n'; Value = 10 },
            @{ Col = 6; Type = 'n'; Value = 3 },
            @{ Col = 7; Type = 'n'; Value = 50 },
            @{ Col = 8; Type = 'n'; Value = 0 },
            @{ Col = 9; Type = 'n'; Value = 0 },
            @{ Col = 10; Type = 'n'; Value = 100 },
            @{ Col = 11; Type = 'n'; Value = 200 },
            @{ Col = 12; Type = 'b'; Value = $false },
            @{ Col = 13; Type = 'n'; Value = -2 },
            @{ Col = 14; Type = 'n'; Value = 100 },
            @{ Col = 15; Type = 'n'; Value = 256 },
            @{ Col = 16; Type = 'n'; Value = 10 },
            @{ Col = 17; Type = 's'; Value = "<-parameter / result->" },
            @{ Col = 18; Type = 'n'; Value = 9 },
            @{ Col = 19; Type = 'n'; Value = -1; Fmt = '0.000000' },
            @{ Col = 20; Type = 'n'; Value = 51.25487888560576; Fmt = '0.000000' },
            @{ Col = 21; Type = 'n'; Value = 17.52092043558757; Fmt = '0.000000' },
            @{ Col = 22; Type = 'n'; Value = 0.1576882839202881; Fmt = '0.000000' },
            @{ Col = 23; Type = 'n'; Value = 0.0001018047332763672; Fmt = '0.000000' },
            @{ Col = 24; Type = 'n'; Value = 0.101953125; Fmt = '0.000000' },
            @{ Col = 25; Type = 'n'; Value = 0.05097656249999999; Fmt = '0.000000' },
            @{ Col = 26; Type = 'n'; Value = 0; Fmt = '0.000000' },
            @{ Col = 27; Type = 'n'; Value = 1.225343465805054; Fmt = '0.000' },
            @{ Col = 28; Type = 'n'; Value = 1044; Fmt = '0.00' },
            @{ Col = 29; Type = 'n'; Value = 0; Fmt = '0.00' }
        )
    },
    @{
        Row = 301
        Cells = @(
            @{ Col = 1; Type = 'n'; Value = 2 },
            @{ Col = 2; Type = 'n'; Value = 1 },
            @{ Col = 3; Type = 'n'; Value = 2 },
            @{ Col = 4; Type = 'n'; Value = 10 },
            @{ Col = 5; Type = 'n'; Value = 10 },
            @{ Col = 6; Type = 'n'; Value = 3 },
            @{ Col = 7; Type = 'n'; Value = 50 },
            @{ Col = 8; Type = 'n'; Value = 0 },
            @{ Col = 9; Type = 'n'; Value = 0 },
            @{ Col = 10; Type = 'n'; Value = 100 },
            @{ Col = 11; Type = 'n'; Value = 200 },
            @{ Col = 12; Type = 'b'; Value = $false },
            @{ Col = 13; Type = 'n'; Value = -2 },
            @{ Col = 14; Type = 'n'; Value = 100 },
            @{ Col = 15; Type = 'n'; Value = 256 },
            @{ Col = 16; Type = 'n'; Value = 10 },
            @{ Col = 17; Type = 's'; Value = "<-parameter / result->" },
            @{ Col = 18; Type = 'n'; Value = 9 },
            @{ Col = 19; Type = 'n'; Value = -1; Fmt = '0.000000' },
            @{ Col = 20; Type = 'n'; Value = 51.34789298562442; Fmt = '0.000000' },
            @{ Col = 21; Type = 'n'; Value = 17.51774152119955; Fmt = '0.000000' },
            @{ Col = 22; Type = 'n'; Value = 0.1576596736907959; Fmt = '0.000000' },
            @{ Col = 23; Type = 'n'; Value = [double]"9.059906005859375e-05"; Fmt = '0.000000' },
            @{ Col = 24; Type = 'n'; Value = 0.101953125; Fmt = '0.000000' },
            @{ Col = 25; Type = 'n'; Value = 0.05097656249999999; Fmt = '0.000000' },
            @{ Col = 26; Type = 'n'; Value = 0; Fmt = '0.000000' },
            @{ Col = 27; Type = 'n'; Value = 1.224569797515869; Fmt = '0.000' },
            @{ Col = 28; Type = 'n'; Value = 1044; Fmt = '0.00' },
            @{ Col = 29; Type = 'n'; Value = 0; Fmt = '0.00' }
        )
    },
    @{
        Row = 302
        Cells = @(
            @{ Col = 1; Type = 'n'; Value = 2 },
            @{ Col = 2; Type = 'n'; Value = 1 },
            @{ Col = 3; Type = 'n'; Value = 2 },
            @{ Col = 4; Type = 'n'; Value = 10 },
            @{ Col = 5; Type = 'n'; Value = 10 },
            @{ Col = 6; Type = 'n'; Value = 3 },
            @{ Col = 7; Type = 'n'; Value = 50 },
            @{ Col = 8; Type = 'n'; Value = 0 },
            @{ Col = 9; Type = 'n'; Value = 0 },
            @{ Col = 10; Type = 'n'; Value = 100 },
            @{ Col = 11; Type = 'n'; Value = 200 },
            @{ Col = 12; Type = 'b'; Value = $false },
            @{ Col = 13; Type = 'n'; Value = -2 },
            @{ Col = 14; Type = 'n'; Value = 100 },
            @{ Col = 15; Type = 'n'; Value = 256 },
            @{ Col = 16; Type = 'n'; Value = 10 },
            @{ Col = 17; Type = 's'; Value = "<-parameter / result->" },
            @{ Col = 18; Type = 'n'; Value = 7 },
            @{ Col = 19; Type = 'n'; Value = -1; Fmt = '0.000000' },
            @{ Col = 20; Type = 'n'; Value = 31.04332923889159; Fmt = '0.000000' },
            @{ Col = 21; Type = 'n'; Value = 10.54476601736886; Fmt = '0.000000' },
            @{ Col = 22; Type = 'n'; Value = 0.07381336212158202; Fmt = '0.000000' },
            @{ Col = 23; Type = 'n'; Value = [double]"7.486343383789062e-05"; Fmt = '0.000000' },
            @{ Col = 24; Type = 'n'; Value = 0.047578125; Fmt = '0.000000' },
            @{ Col = 25; Type = 'n'; Value = 0.0237890625; Fmt = '0.000000' },
            @{ Col = 26; Type = 'n'; Value = 0; Fmt = '0.000000' },
            @{ Col = 27; Type = 'n'; Value = 1.221666812896729; Fmt = '0.000' },
            @{ Col = 28; Type = 'n'; Value = 487.2; Fmt = '0.00' },
            @{ Col = 29; Type = 'n'; Value = 0; Fmt = '0.00' }
        )
    }
)

foreach ($row in $newRows) {
    foreach ($cell in $row.Cells) {
        $target = $ws.Cells.Item($row.Row, $cell.Col)
        $target.Value2 = $cell.Value
        if ($cell.ContainsKey('Fmt')) {
            $target.NumberFormat = $cell.Fmt
        }
    }
}
